$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# PMT sheet: selection change only (A2:G8 -> C36)
# ------------------------------------------------------------------
$wsPMT = $wb.Worksheets.Item("PMT")
$wsPMT.Activate()
$wsPMT.Range("C36").Select()

# ------------------------------------------------------------------
# CALC sheet: rewrite column-B "copy-down" formulas so each row
# references the originating sheet's row-2 / row-3 (or PMT's payment
# rows) directly instead of chaining off column A. This breaks the
# previous B-column shared-formula group (si) on each row; Excel /
# the engine renumbers the remaining shared-formula ids for columns
# A and E automatically on save.
# ------------------------------------------------------------------
$wsCALC = $wb.Worksheets.Item("CALC")

# ORIG block (rows 5-12)
$wsCALC.Range("B5").Formula  = '=IF(A5="","",ORIG!$B$2)'
$wsCALC.Range("B6").Formula  = '=IF(A6="","",ORIG!$B$2)'
$wsCALC.Range("B7").Formula  = '=IF(A7="","",ORIG!$B$2)'
$wsCALC.Range("B8").Formula  = '=IF(A8="","",ORIG!$B$2)'
$wsCALC.Range("B9").Formula  = '=IF(A9="","",ORIG!$B$3)'
$wsCALC.Range("B10").Formula = '=IF(A10="","",ORIG!$B$3)'
$wsCALC.Range("B11").Formula = '=IF(A11="","",ORIG!$B$3)'
$wsCALC.Range("B12").Formula = '=IF(A12="","",ORIG!$B$3)'

# INT_ACC block (rows 17-24)
$wsCALC.Range("B17").Formula = '=IF(A17="","",INT_ACC!$B$2)'
$wsCALC.Range("B18").Formula = '=IF(A18="","",INT_ACC!$B$2)'
$wsCALC.Range("B19").Formula = '=IF(A19="","",INT_ACC!$B$2)'
$wsCALC.Range("B20").Formula = '=IF(A20="","",INT_ACC!$B$2)'
$wsCALC.Range("B21").Formula = '=IF(A21="","",INT_ACC!$B$3)'
$wsCALC.Range("B22").Formula = '=IF(A22="","",INT_ACC!$B$3)'
$wsCALC.Range("B23").Formula = '=IF(A23="","",INT_ACC!$B$3)'
$wsCALC.Range("B24").Formula = '=IF(A24="","",INT_ACC!$B$3)'

# PMT block (rows 29-36)
$wsCALC.Range("B29").Formula = '=IF(A29="","",PMT!B4)'
$wsCALC.Range("B30").Formula = '=IF(A30="","",PMT!B8)'
$wsCALC.Range("B31").Formula = '=IF(A31="","",PMT!B3)'
$wsCALC.Range("B32").Formula = '=IF(A32="","",PMT!B7)'
$wsCALC.Range("B33").Formula = '=IF(A33="","",PMT!B2)'
$wsCALC.Range("B34").Formula = '=IF(A34="","",PMT!B3)'
$wsCALC.Range("B35").Formula = '=IF(A35="","",PMT!B5)'
$wsCALC.Range("B36").Formula = '=IF(A36="","",PMT!B6)'

# RC block (rows 41-54) -- 'RC' needs doubled single-quotes inside the
# single-quoted PowerShell literal to produce a literal apostrophe.
$wsCALC.Range("B41").Formula = '=IF(A41="","",''RC''!$B$2)'
$wsCALC.Range("B42").Formula = '=IF(A42="","",''RC''!$B$2)'
$wsCALC.Range("B43").Formula = '=IF(A43="","",''RC''!$B$2)'
$wsCALC.Range("B44").Formula = '=IF(A44="","",''RC''!$B$2)'
$wsCALC.Range("B45").Formula = '=IF(A45="","",''RC''!$B$2)'
$wsCALC.Range("B46").Formula = '=IF(A46="","",''RC''!$B$2)'
$wsCALC.Range("B47").Formula = '=IF(A47="","",''RC''!$B$2)'
$wsCALC.Range("B48").Formula = '=IF(A48="","",''RC''!$B$3)'
$wsCALC.Range("B49").Formula = '=IF(A49="","",''RC''!$B$3)'
$wsCALC.Range("B50").Formula = '=IF(A50="","",''RC''!$B$3)'
$wsCALC.Range("B51").Formula = '=IF(A51="","",''RC''!$B$3)'
$wsCALC.Range("B52").Formula = '=IF(A52="","",''RC''!$B$3)'
$wsCALC.Range("B53").Formula = '=IF(A53="","",''RC''!$B$3)'
$wsCALC.Range("B54").Formula = '=IF(A54="","",''RC''!$B$3)'

# CALC sheet view: tabSelected moves away from here (Expected1 becomes
# the active sheet below); update scroll position + selection.
$wsCALC.Activate()
$wsCALC.Range("C22").Select()
$excel.ActiveWindow.ScrollRow = 15

# ------------------------------------------------------------------
# o_Transaction sheet: selection + scroll position change only
# ------------------------------------------------------------------
$wsTxn = $wb.Worksheets.Item("o_Transaction")
$wsTxn.Activate()
$wsTxn.Range("D42").Select()
$excel.ActiveWindow.ScrollRow = 9

# ------------------------------------------------------------------
# Expected1 sheet: adjust a handful of transaction dates/amounts, plus
# selection + scroll/active-tab change (this sheet becomes active).
# ------------------------------------------------------------------
$wsExp = $wb.Worksheets.Item("Expected1")
$wsExp.Activate()

$wsExp.Range("B42").Value = 44592
$wsExp.Range("E42").Value = 1000
$wsExp.Range("B43").Value = 44592
$wsExp.Range("B44").Value = 44593
$wsExp.Range("B45").Value = 44593
$wsExp.Range("B47").Value = 44593
$wsExp.Range("B49").Value = 44594
$wsExp.Range("B50").Value = 44593
$wsExp.Range("B51").Value = 44593
$wsExp.Range("B52").Value = 44593
$wsExp.Range("B53").Value = 44593
$wsExp.Range("B54").Value = 44593
$wsExp.Range("B55").Value = 44593
$wsExp.Range("B56").Value = 44593
$wsExp.Range("B57").Value = 44594
$wsExp.Range("B58").Value = 44594
$wsExp.Range("B59").Value = 44594
$wsExp.Range("B60").Value = 44594
$wsExp.Range("B61").Value = 44594
$wsExp.Range("B62").Value = 44594
$wsExp.Range("B63").Value = 44594
$wsExp.Range("B83").Value = 44593
$wsExp.Range("B84").Value = 44594

$wsExp.Range("F57").Select()

$wb.Save()
